{"js": "// Remove the trailing footer paragraphs from the document body:\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n//   - the (now orphaned) empty paragraph that sat between the\n//     \"LOB1037: ...\" requirement line and the \"Ver no Jupiter ...\" line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targets = [];\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const text = (items[i].text || \"\").trim();\n  if (text === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    targets.push(items[i]);\n    // the blank paragraph immediately preceding this one is the spacer\n    // paragraph that should be removed along with the footer block\n    if (i - 1 >= 0 && (items[i - 1].text || \"\").trim() === \"\") {\n      targets.push(items[i - 1]);\n    }\n  } else if (text.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    targets.push(items[i]);\n  }\n}\n\nfor (const p of targets) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing footer paragraphs from the document body:\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n#   - the (now orphaned) empty paragraph that sat between the\n#     \"LOB1037: ...\" requirement line and the \"Ver no Jupiter ...\" line.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n\n    if ($t -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $p.Range.Delete()\n\n        # the blank paragraph immediately preceding this one is the spacer\n        # paragraph that should be removed along with the footer block\n        if ($i - 1 -ge 1) {\n            $prev = $d.Paragraphs.Item($i - 1)\n            if ($prev.Range.Text.Trim() -eq \"\") {\n                $prev.Range.Delete()\n            }\n        }\n    } elseif ($t.Contains(\"Contact: luizeleno@usp.br\")) {\n        $p.Range.Delete()\n    }\n}\n"}
